$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stundenerfassung")

# Append two new rows (105, 106) below the existing data, inheriting the
# formatting of the row directly above (matches Excel's default
# "insert copies format from above" behaviour).
$ws.Rows.Item(105).Insert(-4121, -4163) | Out-Null
$ws.Rows.Item(106).Insert(-4121, -4163) | Out-Null

$ws.Cells.Item(105, 1).Value = 42959
$ws.Cells.Item(105, 2).Value = "ETIC2"
$ws.Cells.Item(105, 3).Value = "Design View Model"
$ws.Cells.Item(105, 4).Value = 5

$ws.Cells.Item(106, 1).Value = 42959
$ws.Cells.Item(106, 2).Value = "ETIC2"
$ws.Cells.Item(106, 3).Value = "Anbindung SoftwareVersionsDatabase"
$ws.Cells.Item(106, 4).Value = 2

# Refresh the view / database context: move the active selection from
# B109 to C109.
$ws.Range("C109").Select() | Out-Null
